# Add language code information to the help
# The underlying Info.docx / Info.pptx help files grew (extra language-code
# documentation), so the character counts billed for translating them went
# up. Reflect the new counts on the Info.xlsx summary sheet and leave the
# view positioned where the edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "characters billed" counts for the Word and PowerPoint help files.
$ws.Range("D3").Value = 1035
$ws.Range("D4").Value = 1034

# Column A / B were auto-sized slightly differently after the edit.
$ws.Columns.Item(1).ColumnWidth = 12.59
$ws.Columns.Item(2).ColumnWidth = 15.25

# Leave the cursor/selection where the author made the change and scroll
# the view over a column.
$ws.Range("D5").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
